$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.280843734741211
$ws.Range("B1").Value = 2.28861141204834
$ws.Range("C1").Value = 4.704154491424561
$ws.Range("D1").Value = 3.029700040817261
$ws.Range("E1").Value = 1.342018365859985
